$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 30: add date, change job no. from 4 to 1
$ws.Range("A30").Value = 43914
$ws.Range("A30").NumberFormat = "m/d;@"
$ws.Range("B30").Value = "1"
$ws.Range("B30").NumberFormat = "@"

# Row 31: change job no. from 5 to 2
$ws.Range("B31").Value = "2"
$ws.Range("B31").NumberFormat = "@"

# Row 32: new entry - "Draw the nearest path"
$ws.Range("B32").Value = "3"
$ws.Range("B32").NumberFormat = "@"
$ws.Range("C32").Value = "Draw the nearest path"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("D32").Value = 0.79861111111111116
$ws.Range("D32").NumberFormat = "[$-409]h:mm\ AM/PM;@"
$ws.Range("E32").Value = 0.84027777777777779
$ws.Range("E32").NumberFormat = "[$-409]h:mm\ AM/PM;@"
$ws.Range("F32").Value = 1
$ws.Range("F32").NumberFormat = "0.00"

# Update selection to match the new state
$ws.Range("F16:F26").Select()
